$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 200
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()

$ws.Range("H64").Value = 4418
$ws.Range("J64").Value = 4866.5
$ws.Range("L64").Value = 4866.5
$ws.Range("N64").Value = -5362.5

$ws.Range("H67").Value = 4418
$ws.Range("J67").Value = 4866.5
$ws.Range("L67").Value = 4866.5
$ws.Range("N67").Value = -6582.5

$ws.Range("H74").Value = 2440.6667
$ws.Range("I74").Value = 1694.6111
$ws.Range("K74").Value = 1694.6111
$ws.Range("M74").Value = -758.6111000000001

$ws.Range("H76").Value = 3271.4285
$ws.Range("I76").Value = 3233.3333
$ws.Range("K76").Value = 3233.3333
$ws.Range("M76").Value = -2918.3333

$ws.Range("H77").Value = 2440.6667
$ws.Range("I77").Value = 1694.6111
$ws.Range("K77").Value = 8473.0555
$ws.Range("M77").Value = -3793.0555

$ws.Range("H79").Value = 3271.4285
$ws.Range("I79").Value = 3233.3333
$ws.Range("K79").Value = 3233.3333
$ws.Range("M79").Value = -2141.3333

$ws.Range("H88").Value = 1675
$ws.Range("I88").Value = 800
$ws.Range("J88").Value = 1784.375
$ws.Range("K88").Value = 800
$ws.Range("L88").Value = 1784.375
$ws.Range("M88").Value = -394
$ws.Range("N88").Value = -2596.375

$ws.Range("H91").Value = 1675
$ws.Range("I91").Value = 800
$ws.Range("J91").Value = 1784.375
$ws.Range("K91").Value = 800
$ws.Range("L91").Value = 1784.375
$ws.Range("M91").Value = 604
$ws.Range("N91").Value = -4592.375

$ws.Range("H129").Value = 854.6613
$ws.Range("J129").Value = 870.65515
$ws.Range("L129").Value = 2611.96545
$ws.Range("N129").Value = -12611.96545

$ws.Range("H138").Value = 3045.8572
$ws.Range("I138").Value = 846.75
$ws.Range("J138").Value = 3329.6128
$ws.Range("K138").Value = 2540.25
$ws.Range("L138").Value = 9988.838400000001
$ws.Range("M138").Value = 2599.75
$ws.Range("N138").Value = -20268.8384

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 30174.525
$ws.Range("I32").Value = 38289.965
$ws.Range("J32").Value = 7451.3
$ws.Range("K32").Value = 38289.965
$ws.Range("L32").Value = 7451.3
$ws.Range("M32").Value = -38002.965
$ws.Range("N32").Value = -8025.3

$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").ClearContents()
$ws.Range("N36").Value = 0

$ws.Range("H82").Value = 28181
$ws.Range("J82").Value = 28181
$ws.Range("L82").Value = 28181
$ws.Range("N82").Value = -28903

$ws.Range("H85").Value = 28181
$ws.Range("J85").Value = 28181
$ws.Range("L85").Value = 28181
$ws.Range("N85").Value = -30677

$ws.Range("H102").Value = 3440.5715
$ws.Range("I102").Value = 2253.1428
$ws.Range("K102").Value = 2253.1428
$ws.Range("M102").Value = -631.1428000000001

$ws.Range("H110").Value = 1699.1
$ws.Range("I110").Value = 1438.6154
$ws.Range("J110").Value = 2182.8572
$ws.Range("K110").Value = 1438.6154
$ws.Range("L110").Value = 2182.8572
$ws.Range("M110").Value = 606.3846000000001
$ws.Range("N110").Value = -6272.8572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 29383.111
$ws.Range("I86").Value = 46672.727
$ws.Range("J86").Value = 2213.7144
$ws.Range("K86").Value = 46672.727
$ws.Range("L86").Value = 2213.7144
$ws.Range("M86").Value = -45549.727
$ws.Range("N86").Value = -4459.7144

$ws.Range("H89").Value = 29383.111
$ws.Range("I89").Value = 46672.727
$ws.Range("J89").Value = 2213.7144
$ws.Range("K89").Value = 233363.635
$ws.Range("L89").Value = 11068.572
$ws.Range("M89").Value = -227747.635
$ws.Range("N89").Value = -22300.572

$ws.Range("H107").Value = 646.6070999999999
$ws.Range("I107").Value = 633.8
$ws.Range("K107").Value = 633.8
$ws.Range("M107").Value = 1286.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 221.14285
$ws.Range("I7").Value = 69.59999999999999
$ws.Range("J7").Value = 600
$ws.Range("K7").Value = 69.59999999999999
$ws.Range("L7").Value = 600
$ws.Range("M7").Value = 43.40000000000001
$ws.Range("N7").Value = -826

$ws.Range("H31").Value = 9403.130999999999
$ws.Range("I31").Value = 19071.666
$ws.Range("J31").Value = 3187.6428
$ws.Range("K31").Value = 19071.666
$ws.Range("L31").Value = 3187.6428
$ws.Range("M31").Value = -18776.666
$ws.Range("N31").Value = -3777.6428

$ws.Range("H34").Value = 9403.130999999999
$ws.Range("I34").Value = 19071.666
$ws.Range("J34").Value = 3187.6428
$ws.Range("K34").Value = 19071.666
$ws.Range("L34").Value = 3187.6428
$ws.Range("M34").Value = -18869.666
$ws.Range("N34").Value = -3591.6428

$ws.Range("H41").Value = 22666.334
$ws.Range("J41").Value = 26499.5
$ws.Range("L41").Value = 26499.5
$ws.Range("N41").Value = -27355.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 6000040.5
$ws.Range("I4").Value = 67.333336
$ws.Range("K4").Value = 202.000008
$ws.Range("M4").Value = -90.00000800000001

$ws.Range("H12").Value = 90.94444
$ws.Range("J12").Value = 121.07692
$ws.Range("L12").Value = 363.23076
$ws.Range("N12").Value = -709.23076

$ws.Range("H40").Value = 547.1429000000001
$ws.Range("I40").Value = 132.5
$ws.Range("J40").Value = 1100
$ws.Range("K40").Value = 530
$ws.Range("L40").Value = 4400
$ws.Range("M40").Value = -461
$ws.Range("N40").Value = -4538

$ws.Range("H68").Value = 1345.1904
$ws.Range("J68").Value = 1415.7894
$ws.Range("L68").Value = 4247.3682
$ws.Range("N68").Value = -5869.3682

$ws.Range("H69").Value = 1666.6666
$ws.Range("J69").Value = 1750
$ws.Range("L69").Value = 5250
$ws.Range("N69").Value = -6872

$ws.Range("H71").Value = 1345.1904
$ws.Range("J71").Value = 1415.7894
$ws.Range("L71").Value = 12742.1046
$ws.Range("N71").Value = -20854.1046

$ws.Range("H72").Value = 1666.6666
$ws.Range("J72").Value = 1750
$ws.Range("L72").Value = 15750
$ws.Range("N72").Value = -23862

$ws.Range("H131").Value = 164780.48
$ws.Range("J131").Value = 179405.64
$ws.Range("L131").Value = 538216.92
$ws.Range("N131").Value = -548296.92

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 13098.637
$ws.Range("J80").Value = 5316
$ws.Range("L80").Value = 5316
$ws.Range("N80").Value = -7312

$ws.Range("H83").Value = 13098.637
$ws.Range("J83").Value = 5316
$ws.Range("L83").Value = 26580
$ws.Range("N83").Value = -36564

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").ClearContents()
$ws.Range("N130").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 9900
$ws.Range("I45").Value = 1800
$ws.Range("J45").Value = 18000
$ws.Range("K45").Value = 1800
$ws.Range("L45").Value = 18000
$ws.Range("M45").Value = -1393
$ws.Range("N45").Value = -18814

$ws.Range("H46").Value = 1295.3572
$ws.Range("I46").Value = 833.5
$ws.Range("J46").Value = 2450
$ws.Range("K46").Value = 833.5
$ws.Range("L46").Value = 2450
$ws.Range("M46").Value = -645.5
$ws.Range("N46").Value = -2826

$ws.Range("H82").Value = 3699
$ws.Range("I82").Value = 4373.75
$ws.Range("K82").Value = 4373.75
$ws.Range("M82").Value = -4012.75

$ws.Range("H85").Value = 3699
$ws.Range("I85").Value = 4373.75
$ws.Range("K85").Value = 4373.75
$ws.Range("M85").Value = -3125.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1354.4222
$ws.Range("I132").Value = 1070.0303
$ws.Range("J132").Value = 2136.5
$ws.Range("K132").Value = 3210.0909
$ws.Range("L132").Value = 6409.5
$ws.Range("M132").Value = -680.0908999999997
$ws.Range("N132").Value = -11469.5
